$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the old per-occurrence "Quiz N" reading entries in column E.
# (Quiz 1 moves from row 5 down to row 10; Quiz 2-7 are dropped entirely.)
$ws.Range("E5").Clear()
$ws.Range("E8").Clear()
$ws.Range("E11").Clear()
$ws.Range("E14").Clear()
$ws.Range("E17").Clear()
$ws.Range("E20").Clear()
$ws.Range("E23").Clear()

# Row 10 (Wed 03 Sep 2014): add the "Quiz 1" reading note and swap the
# reading from "LBNL Reading" to the new Meier, 1984 reading.
$ws.Range("E10").Value = "Quiz 1"
$ws.Range("G10").Value = "Meier, 1984"

# Update view state to match: no frozen/scrolled topLeftCell, selection on G10.
$ws.Range("G10").Select()
